# VMC and Prime Video Test Cases
# Adds a new test-case row (row 33) to Sheet1 with:
#   Role=Student, Board=CBSE, Class Name=12, Role Type=Single,
#   Subscription=Guru, Mobile=Yes, Username=parul.s, Password=123456

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A33").Value = "Student"
$ws.Range("B33").Value = "CBSE"

$ws.Range("C33").Value = 12
$ws.Range("C33").NumberFormat = "@"

$ws.Range("D33").Value = "Single"
$ws.Range("E33").Value = "Guru"
$ws.Range("F33").Value = "Yes"
$ws.Range("G33").Value = "parul.s"

$ws.Range("H33").Value = 123456
$ws.Range("H33").NumberFormat = "@"

# Move the view/selection down to the newly added row, mirroring the
# author's on-screen position when the change was saved.
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$win.ScrollColumn = 1
$ws.Range("D33").Select() | Out-Null
